$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B ("Buying Opportunity") values for rows 2-18 ---
$ws.Range("B2").Value  = "NSE:AVROIND"
$ws.Range("B3").Value  = "NSE:AXISBNKETF"
$ws.Range("B4").Value  = "NSE:DEVYANI"
$ws.Range("B5").Value  = "NSE:DRCSYSTEMS"
$ws.Range("B6").Value  = "NSE:EMKAY"
$ws.Range("B7").Value  = "NSE:GLENMARK"
$ws.Range("B8").Value  = "NSE:HEG"
$ws.Range("B9").Value  = "NSE:HINDZINC"
$ws.Range("B10").Value = "NSE:KANORICHEM"
$ws.Range("B11").Value = "NSE:KAUSHALYA"
$ws.Range("B12").Value = "NSE:MAGADSUGAR"
$ws.Range("B13").Value = "NSE:MAHEPC"
$ws.Range("B14").Value = "NSE:MATRIMONY"
$ws.Range("B15").Value = "NSE:MINDTECK"
$ws.Range("B16").Value = "NSE:MODIRUBBER"
$ws.Range("B17").Value = "NSE:NEWGEN"
$ws.Range("B18").Value = "NSE:NUVAMA"

# --- Update column C ("support Zone") values for rows 3-18 (row 2 is unchanged) ---
$ws.Range("C3").Value  = "NSE:ASALCBR"
$ws.Range("C4").Value  = "NSE:ATGL"
$ws.Range("C5").Value  = "NSE:ATUL"
$ws.Range("C6").Value  = "NSE:AYMSYNTEX"
$ws.Range("C7").Value  = "NSE:CENTENKA"
$ws.Range("C8").Value  = "NSE:DPSCLTD"
$ws.Range("C9").Value  = "NSE:DYCL"
$ws.Range("C10").Value = "NSE:ELGIEQUIP"
$ws.Range("C11").Value = "NSE:GMMPFAUDLR"
$ws.Range("C12").Value = "NSE:HPL"
$ws.Range("C13").Value = "NSE:INDOBORAX"
$ws.Range("C14").Value = "NSE:LXCHEM"
$ws.Range("C15").Value = "NSE:ONMOBILE"
$ws.Range("C16").Value = "NSE:PATELENG"
$ws.Range("C17").Value = "NSE:PRIMESECU"
$ws.Range("C18").Value = "NSE:PTL"

# --- New rows 19-21: only column A (index) and column C (support Zone) are populated ---
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19

$ws.Range("C19").Value = "NSE:QUESS"
$ws.Range("C20").Value = "NSE:RAYMOND"
$ws.Range("C21").Value = "NSE:RICOAUTO"

# Copy the formatting (bold, centered, bordered "index" style) from A18 down to the new index cells
$ws.Range("A18").Copy()
$ws.Range("A19:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Touch the remaining (empty) cells in the new rows so they materialize like the existing sparse rows
$ws.Cells.Item(19,2).Font.Bold = $false
$ws.Cells.Item(19,4).Font.Bold = $false
$ws.Cells.Item(19,5).Font.Bold = $false
$ws.Cells.Item(19,6).Font.Bold = $false

$ws.Cells.Item(20,2).Font.Bold = $false
$ws.Cells.Item(20,4).Font.Bold = $false
$ws.Cells.Item(20,5).Font.Bold = $false
$ws.Cells.Item(20,6).Font.Bold = $false

$ws.Cells.Item(21,2).Font.Bold = $false
$ws.Cells.Item(21,4).Font.Bold = $false
$ws.Cells.Item(21,5).Font.Bold = $false
$ws.Cells.Item(21,6).Font.Bold = $false
